# RPA datasets push 2023-10-19
# Insert a new IPO record (퓨릿, listed 2023-10-18) at the top of the data
# table (row 2), pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2..24) down to (3..25) by inserting a fresh row.
$ws.Rows.Item(2).Insert()

# Populate the new row with the new listing's data.
$ws.Range("A2").Value = "'2023-10-18"
$ws.Range("B2").Value = "퓨릿"
$ws.Range("C2").Value = "코스닥"
$ws.Range("D2").Value = 442.659
$ws.Range("E2").Value = "미래"
$ws.Range("F2").Value = 442.659
$ws.Range("G2").Value = "-"
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "-"
$ws.Range("K2").Value = "대표"
$ws.Range("L2").Value = "-"
$ws.Range("M2").Value = 10700
$ws.Range("N2").Value = 100
$ws.Range("O2").Value = "'2023-10-05"
$ws.Range("P2").Value = "'2023-10-11"
$ws.Range("Q2").Value = 3102750

# Drop the formatting that leaked in from the Insert() (quote-prefix /
# number-format styles) so the new row matches the plain, unstyled data
# rows used throughout the rest of the table.
$ws.Range("A2:Q2").ClearFormats()
